# fix reg detail view if in individual mode
#
# 1. Swap the "id_type" (select_one) question and the "id_number" (string)
#    question on the "survey" sheet: row 10 <-> row 11 (columns C, D, E, H).
# 2. Make the "survey" sheet the active tab (it was "model"), scrolled so
#    row 18 is at the top, with C28 selected.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")

# --- swap row 10 and row 11 contents (type / name / prompt / values_list) ---
$c10 = $survey.Range("C10").Text
$d10 = $survey.Range("D10").Text
$e10 = $survey.Range("E10").Text
$h10 = $survey.Range("H10").Text

$c11 = $survey.Range("C11").Text
$d11 = $survey.Range("D11").Text
$e11 = $survey.Range("E11").Text
$h11 = $survey.Range("H11").Text

$survey.Range("C10").Value = $c11
$survey.Range("D10").Value = $d11
$survey.Range("E10").Value = $e11
$survey.Range("H10").Value = $h11

$survey.Range("C11").Value = $c10
$survey.Range("D11").Value = $d10
$survey.Range("E11").Value = $e10
$survey.Range("H11").Value = $h10

# --- activate the survey sheet/view (previously "model" was active) ---
$survey.Activate()
$survey.Range("C28").Select()
